# Added degree plans in SpreadSheet/Sample Data.xlsx
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# DegreePlan sheet: two new degree-plan rows for student 533982
# ----------------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("DegreePlan")

# Write the new plan-name / abbreviation text column-first so the shared
# string table fills in the same order the workbook was authored in.
$wsPlan.Cells.Item(4, 3).Value = "Easy study plan"
$wsPlan.Cells.Item(5, 3).Value = "1 year plan"
$wsPlan.Cells.Item(4, 4).Value = "Study with break"
$wsPlan.Cells.Item(5, 4).Value = "Complete in one year"

$wsPlan.Cells.Item(4, 1).Value = 7253
$wsPlan.Cells.Item(4, 2).Value = 533982
$wsPlan.Cells.Item(4, 5).Value = 3
$wsPlan.Rows.Item(4).RowHeight = 15

$wsPlan.Cells.Item(5, 1).Value = 7254
$wsPlan.Cells.Item(5, 2).Value = 533982
$wsPlan.Cells.Item(5, 5).Value = 3
$wsPlan.Rows.Item(5).RowHeight = 15

# DegreePlanName column got wider once the longer plan names were added
$wsPlan.Columns.Item(4).ColumnWidth = 19.5

# ----------------------------------------------------------------------
# Slot sheet: fill in the slot schedule for the three new degree plans
# (student 7252's remaining slots plus all of 7253 and 7254)
# ----------------------------------------------------------------------
$wsSlot = $wb.Worksheets.Item("Slot")

$slotData = @(
    ,@(13, 12, 7252, 1, 664, "p", 15)
    ,@(14, 13, 7252, 1, 64, "p", 15)
    ,@(15, 14, 7252, 1, 10, "p", 15)
    ,@(16, 15, 7252, 2, 691, "p", 15)
    ,@(17, 16, 7252, 2, 555, "p", 15)
    ,@(18, 17, 7252, 2, 618, "p", 15)
    ,@(19, 18, 7252, 3, 460, "A", 15)
    ,@(20, 19, 7252, 3, 542, "A", 15)
    ,@(21, 20, 7252, 3, 563, "A", 15.75)
    ,@(22, 21, 7252, 4, 560, "C", 15.75)
    ,@(23, 22, 7252, 4, 20, "C", 15.75)
    ,@(24, 23, 7252, 4, 692, "C", 15.75)
    ,@(25, 24, 7253, 1, 542, "p", 15.75)
    ,@(26, 25, 7253, 1, 563, "P", 15.75)
    ,@(27, 26, 7253, 1, 460, "P", 15.75)
    ,@(28, 27, 7253, 2, 560, "A", 15.75)
    ,@(29, 28, 7253, 2, 664, "A", 15.75)
    ,@(30, 29, 7253, 2, 64, "A", 15.75)
    ,@(31, 30, 7253, 3, 691, "A", 15.75)
    ,@(32, 31, 7253, 3, 10, "A", 15.75)
    ,@(33, 32, 7253, 3, 555, "A", 15.75)
    ,@(34, 33, 7253, 4, 692, "C", 15.75)
    ,@(35, 34, 7253, 4, 20, "C", 15.75)
    ,@(36, 35, 7253, 4, 356, "C", 15.75)
    ,@(37, 36, 7254, 1, 460, "P", 15.75)
    ,@(38, 37, 7254, 1, 542, "P", 15.75)
    ,@(39, 38, 7254, 1, 356, "P", 15.75)
    ,@(40, 39, 7254, 2, 563, "C", 15.75)
    ,@(41, 40, 7254, 2, 560, "C", 15.75)
    ,@(42, 41, 7254, 2, 664, "C", 15.75)
    ,@(43, 42, 7254, 3, 618, "A", 15.75)
    ,@(44, 43, 7254, 3, 555, "A", 15.75)
    ,@(45, 44, 7254, 3, 691, "A", 15.75)
    ,@(46, 45, 7254, 4, 692, "P", 15.75)
    ,@(47, 46, 7254, 4, 10, "P", 15.75)
    ,@(48, 47, 7254, 4, 64, "P", 15.75)
)

foreach ($row in $slotData) {
    $r = $row[0]
    $wsSlot.Cells.Item($r, 1).Value = $row[1]
    $wsSlot.Cells.Item($r, 2).Value = $row[2]
    $wsSlot.Cells.Item($r, 3).Value = $row[3]
    $wsSlot.Cells.Item($r, 3).HorizontalAlignment = -4108
    $wsSlot.Cells.Item($r, 4).Value = $row[4]
    $wsSlot.Cells.Item($r, 4).HorizontalAlignment = -4108
    $wsSlot.Cells.Item($r, 5).Value = $row[5]
    $wsSlot.Cells.Item($r, 5).HorizontalAlignment = -4108
    # rows 13-20 are brand new rows (default height); 21-48 reuse the
    # pre-existing placeholder rows, which already carry ht="15.75"
    if ($r -le 20) {
        $wsSlot.Rows.Item($r).RowHeight = $row[6]
    }
}

# ----------------------------------------------------------------------
# View / selection state changes
# ----------------------------------------------------------------------

# Credit: selection moved to A7
$wsCredit = $wb.Worksheets.Item("Credit")
$wsCredit.Range("A7").Select()

# DegreePlan: selection now spans the whole DegreePlanName column, and it
# is no longer the active tab
$wsPlan.Range("D1:D1048576").Select()

# Slot: selection moved to E48 (scrolled down to the new rows)
$wsSlot.Range("E48").Select()

# StudentTerm becomes the active tab
$wsTerm = $wb.Worksheets.Item("StudentTerm")
$wsTerm.Activate()
